$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.140.79"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "'1.656.84"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").Value = "'218.56"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "'0.5247"
$ws.Range("E6").Value = "  -1.57%  "
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("D9").Value = "'0.06292"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").Value = "'20.55"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "'0.07804"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "'4.492"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").Value = "'1.663.95"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "'1.884.37"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "'0.5553"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").Value = "'0.0₅8006"
$ws.Range("E16").Value = "  -2.56%  "
$ws.Range("D17").Value = "'65.04"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "'26.159.81"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").Value = "'4.635"
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("D21").Value = "'195.28"
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("D22").Value = "'10.10"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").Value = "'5.964"
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").Value = "'146.73"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D26").Value = "'0.1205"
$ws.Range("E26").Value = "  -1.93%  "
$ws.Range("D27").Value = "'7.179"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("D29").Value = "'1.496"
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").Value = "'0.05719"
$ws.Range("E30").Value = "  -2.69%  "
$ws.Range("D31").Value = "'1.269"
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").Value = "'3.484"
$ws.Range("E32").Value = "  -3.10%  "
$ws.Range("D33").Value = "'3.342"
$ws.Range("E33").Value = "  +1.91%  "
$ws.Range("D34").Value = "'1.584"
$ws.Range("E34").Value = "  -1.58%  "
$ws.Range("D35").Value = "'2.803"
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D36").Value = "'0.9510"
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("D37").Value = "'2.418"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "'0.5706"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").Value = "'0.01595"
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("D40").Value = "'5.941"
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("D41").Value = "'1.063.49"
$ws.Range("E41").Value = "  +1.36%  "
$ws.Range("D42").Value = "'0.8466"
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").Value = "'103.31"
$ws.Range("E44").Value = "  -0.72%  "
$ws.Range("D45").Value = "'1.794.64"
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("D46").Value = "'57.81"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.008"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").Value = "'0.05356"
$ws.Range("E48").Value = "  +3.80%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.0₈104"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").Value = "'0.4401"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").Value = "'7.976"
$ws.Range("E51").Value = "  -0.81%  "
